# "Add Quizzes to submodule 4, updated slides"
# This particular slide (the stand-alone "Quiz" slide, #20) is removed from
# the deck; PowerPoint shifts every following slide up by one position, so
# the slide that used to be #21 ("References") becomes the new #20.

$p = $ppt.ActivePresentation

# Locate and remove the "Quiz" slide.
$quizSlide = $p.Slides.Item(20)
$quizSlide.Delete()

# The slide that follows ("References") is now slide 20. Its footer
# placeholder still caches the slide-number text from its old position
# ("21"); refresh it to match its new position ("20").
$refSlide = $p.Slides.Item(20)
foreach ($shp in $refSlide.Shapes) {
    if ($shp.HasTextFrame -and $shp.Type -eq 14 -and $shp.PlaceholderFormat.Type -eq 13) {
        $shp.TextFrame.TextRange.Text = "20"
    }
}
